$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the two string values in the table (C8 and B10)
$ws.Range("C8").Value = "dropdownLangauge"
$ws.Range("B10").Value = "xpath"

# Update the active cell / selection to B10
$ws.Range("B10").Select()
